$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.472.53"
$ws.Cells.Item(2, 5).Value = "  -0.21%  "

$ws.Cells.Item(3, 4).Value = "'1.841.25"
$ws.Cells.Item(3, 5).Value = "  -0.26%  "

$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "  +0.17%  "

$ws.Cells.Item(5, 4).Value = "'261.31"
$ws.Cells.Item(5, 5).Value = "  -0.90%  "

$ws.Cells.Item(6, 5).Value = "  +0.20%  "

$ws.Cells.Item(7, 4).Value = "'0.5353"
$ws.Cells.Item(7, 5).Value = "  +2.58%  "

$ws.Cells.Item(8, 4).Value = "'0.3034"
$ws.Cells.Item(8, 5).Value = "  -6.23%  "

$ws.Cells.Item(9, 4).Value = "'0.06886"
$ws.Cells.Item(9, 5).Value = "  +1.26%  "

$ws.Cells.Item(10, 4).Value = "'17.95"
$ws.Cells.Item(10, 5).Value = "  -4.33%  "

$ws.Cells.Item(11, 4).Value = "'0.07589"
$ws.Cells.Item(11, 5).Value = "  -2.55%  "

$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "'1.854.31"
$ws.Cells.Item(12, 5).Value = "  +0.45%  "

$ws.Cells.Item(13, 2).Value = "Polygon"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(13, 4).Value = "'0.7367"
$ws.Cells.Item(13, 5).Value = "  -5.21%  "

$ws.Cells.Item(14, 4).Value = "'89.70"
$ws.Cells.Item(14, 5).Value = "  +1.46%  "

$ws.Cells.Item(15, 4).Value = "'4.986"
$ws.Cells.Item(15, 5).Value = "  -0.63%  "

$ws.Cells.Item(16, 5).Value = "  +0.22%  "

$ws.Cells.Item(17, 4).Value = "'13.95"
$ws.Cells.Item(17, 5).Value = "  -0.10%  "

$ws.Cells.Item(18, 5).Value = "  +0.21%  "

$ws.Cells.Item(19, 4).Value = "'0.000007933"
$ws.Cells.Item(19, 5).Value = "  -0.33%  "

$ws.Cells.Item(20, 4).Value = "'26.501.63"

$ws.Cells.Item(21, 4).Value = "'2.091.67"
$ws.Cells.Item(21, 5).Value = "  +0.35%  "

$ws.Cells.Item(22, 4).Value = "'4.599"
$ws.Cells.Item(22, 5).Value = "  -0.53%  "

$ws.Cells.Item(23, 4).Value = "'5.981"
$ws.Cells.Item(23, 5).Value = "  -0.28%  "

$ws.Cells.Item(24, 4).Value = "'9.299"
$ws.Cells.Item(24, 5).Value = "  -1.54%  "

$ws.Cells.Item(25, 4).Value = "'143.24"
$ws.Cells.Item(25, 5).Value = "  +0.17%  "

$ws.Cells.Item(26, 4).Value = "'2.213"
$ws.Cells.Item(26, 5).Value = "  +2.34%  "

$ws.Cells.Item(27, 5).Value = "  +0.58%  "

$ws.Cells.Item(28, 4).Value = "'16.97"
$ws.Cells.Item(28, 5).Value = "  -0.26%  "

$ws.Cells.Item(29, 4).Value = "'110.75"
$ws.Cells.Item(29, 5).Value = "  -0.99%  "

$ws.Cells.Item(30, 4).Value = "'4.257"
$ws.Cells.Item(30, 5).Value = "  +1.79%  "

$ws.Cells.Item(31, 4).Value = "'0.08807"
$ws.Cells.Item(31, 5).Value = "  +0.76%  "

$ws.Cells.Item(32, 4).Value = "'4.058"
$ws.Cells.Item(32, 5).Value = "  -1.23%  "

$ws.Cells.Item(33, 5).Value = "  -0.62%  "

$ws.Cells.Item(34, 4).Value = "'2.929"
$ws.Cells.Item(34, 5).Value = "  +2.38%  "

$ws.Cells.Item(35, 4).Value = "'0.7261"
$ws.Cells.Item(35, 5).Value = "  +0.60%  "

$ws.Cells.Item(36, 4).Value = "'1.133"
$ws.Cells.Item(36, 5).Value = "  +0.35%  "

$ws.Cells.Item(37, 4).Value = "'3.105"
$ws.Cells.Item(37, 5).Value = "  +0.17%  "

$ws.Cells.Item(38, 4).Value = "'2.295"
$ws.Cells.Item(38, 5).Value = "  +3.96%  "

$ws.Cells.Item(39, 4).Value = "'0.01717"

$ws.Cells.Item(40, 4).Value = "'0.4763"
$ws.Cells.Item(40, 5).Value = "  -1.67%  "

$ws.Cells.Item(41, 4).Value = "'0.9078"
$ws.Cells.Item(41, 5).Value = "  +1.82%  "

$ws.Cells.Item(42, 4).Value = "'107.94"
$ws.Cells.Item(42, 5).Value = "  -2.85%  "

$ws.Cells.Item(43, 4).Value = "'5.872"
$ws.Cells.Item(43, 5).Value = "  -2.69%  "

$ws.Cells.Item(44, 4).Value = "'1.001"
$ws.Cells.Item(44, 5).Value = "  +0.21%  "

$ws.Cells.Item(45, 4).Value = "'7.450"
$ws.Cells.Item(45, 5).Value = "  -2.18%  "

$ws.Cells.Item(46, 4).Value = "'9.020"
$ws.Cells.Item(46, 5).Value = "  -0.70%  "

$ws.Cells.Item(47, 4).Value = "'0.4105"
$ws.Cells.Item(47, 5).Value = "  -2.30%  "

$ws.Cells.Item(48, 4).Value = "'0.1234"
$ws.Cells.Item(48, 5).Value = "  -0.30%  "

$ws.Cells.Item(49, 4).Value = "'34.93"
$ws.Cells.Item(49, 5).Value = "  -0.09%  "

$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).Value = "'0.05798"
$ws.Cells.Item(50, 5).Value = "  -1.57%  "

$ws.Cells.Item(51, 2).Value = "EOS"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(51, 4).Value = "'0.8965"
$ws.Cells.Item(51, 5).Value = "  +1.05%  "
